# HOTCARD sheet - Wave update: populate row 4 with the new entity's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOTCARD")

# Plain text fields: a direct Value assignment keeps the cell's existing
# style/number-format untouched and stores the text as-is.
$ws.Range("B4").Value = "National Bank"
$ws.Range("D4").Value = "FISB"
$ws.Range("E4").Value = "Legacy Core"
$ws.Range("F4").Value = "PaymentsOne Debit"
$ws.Range("G4").Value = "Premium"
$ws.Range("H4").Value = "Dallas, TX"
$ws.Range("J4").Value = "Yes"

# Date-/number-looking fields: assigning these directly makes Excel parse
# them as a date/number and reformat the cell. Instead, enter each one into
# a scratch cell with a leading apostrophe (forcing literal text), copy it,
# and paste-special "Values" into the destination so the destination cell
# keeps its original style/number format but receives the literal text.
$helper = $ws.Range("Z1")

$helper.Value = "'11/26/2025"
$helper.Copy()
$ws.Range("A4").PasteSpecial(-4163)

$helper.Value = "'101010"
$helper.Copy()
$ws.Range("C4").PasteSpecial(-4163)

$helper.Value = "'1900"
$helper.Copy()
$ws.Range("I4").PasteSpecial(-4163)

$helper.Clear()
$excel.CutCopyMode = $false
